$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.582.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.238.17'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '115.13'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '282.92'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +7.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.86%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.614'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.76'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.35'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.884'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.574.42'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.236.87'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.773.68'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.85'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.23'
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.18'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +10.87%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.35'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '231.97'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.30'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.09'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.53%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '40.42'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.08%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.29'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.53'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0900'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.56'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +17.88%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.34%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.66'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0372'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.53%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.19%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.41'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.01%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.47'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.234'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.57'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -8.68%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.51'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.649'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +9.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0988'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.471'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +8.53%  '
